$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 271.17648
$ws.Range("I39").Value = 193.3077
$ws.Range("K39").Value = 579.9231
$ws.Range("M39").Value = -283.9231

$ws.Range("H96").Value = 2355.4443
$ws.Range("I96").Value = 2742.8572
$ws.Range("J96").Value = 999.5
$ws.Range("K96").Value = 8228.571599999999
$ws.Range("L96").Value = 2998.5
$ws.Range("M96").Value = -6855.571599999999
$ws.Range("N96").Value = -5744.5

$ws.Range("H132").Value = 1357.8197
$ws.Range("I132").Value = 1327.2632
$ws.Range("J132").Value = 1793.25
$ws.Range("K132").Value = 3981.7896
$ws.Range("L132").Value = 5379.75
$ws.Range("M132").Value = -1451.7896
$ws.Range("N132").Value = -10439.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2603.1458
$ws.Range("I32").Value = 2603.1458
$ws.Range("K32").Value = 2603.1458
$ws.Range("M32").Value = -2316.1458

$ws.Range("H39").Value = 2995.8333
$ws.Range("J39").Value = 5399.6665
$ws.Range("L39").Value = 5399.6665
$ws.Range("N39").Value = -6439.6665

$ws.Range("H63").Value = 1776.4286
$ws.Range("I63").Value = 1776.4286
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1776.4286
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1090.4286
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 1776.4286
$ws.Range("I66").Value = 1776.4286
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8882.143
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5450.143
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 27781846
$ws.Range("I74").Value = 41669984
$ws.Range("K74").Value = 41669984
$ws.Range("M74").Value = -41669110

$ws.Range("H77").Value = 27781846
$ws.Range("I77").Value = 41669984
$ws.Range("K77").Value = 208349920
$ws.Range("M77").Value = -208345552

$ws.Range("H88").Value = 3173.2104
$ws.Range("I88").Value = 3519.5
$ws.Range("J88").Value = 3013.3845
$ws.Range("K88").Value = 3519.5
$ws.Range("L88").Value = 3013.3845
$ws.Range("M88").Value = -3113.5
$ws.Range("N88").Value = -3825.3845

$ws.Range("H91").Value = 3173.2104
$ws.Range("I91").Value = 3519.5
$ws.Range("J91").Value = 3013.3845
$ws.Range("K91").Value = 3519.5
$ws.Range("L91").Value = 3013.3845
$ws.Range("M91").Value = -2115.5
$ws.Range("N91").Value = -5821.3845

$ws.Range("H132").Value = 7992.303
$ws.Range("I132").Value = 4365.091
$ws.Range("K132").Value = 13095.273
$ws.Range("M132").Value = -10565.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 771.5
$ws.Range("I22").Value = 771.5
$ws.Range("K22").Value = 771.5
$ws.Range("M22").Value = -598.5

$ws.Range("H132").Value = 65330.168
$ws.Range("J132").Value = 65330.168
$ws.Range("L132").Value = 65330.168
$ws.Range("N132").Value = -75450.16800000001

$ws.Range("H134").Value = 2647.7646
$ws.Range("I134").Value = 1661.4615
$ws.Range("K134").Value = 4984.3845
$ws.Range("M134").Value = -2449.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 6681.769
$ws.Range("I35").Value = 6681.769
$ws.Range("K35").Value = 6681.769
$ws.Range("M35").Value = -6387.769

$ws.Range("H39").Value = 12500
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9609

$ws.Range("H49").Value = 12500
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9818

$ws.Range("H86").Value = 12498.444
$ws.Range("I86").Value = 5828
$ws.Range("K86").Value = 5828
$ws.Range("M86").Value = -4705

$ws.Range("H89").Value = 12498.444
$ws.Range("I89").Value = 5828
$ws.Range("K89").Value = 29140
$ws.Range("M89").Value = -23524

$ws.Range("H122").Value = 42300.15
$ws.Range("I122").Value = 58266.832
$ws.Range("K122").Value = 174800.496
$ws.Range("M122").Value = -172350.496

$ws.Range("H134").Value = 4762.0625
$ws.Range("I134").Value = 2993.6667
$ws.Range("J134").Value = 8138.091
$ws.Range("K134").Value = 8981.000100000001
$ws.Range("L134").Value = 24414.273
$ws.Range("M134").Value = -6446.000100000001
$ws.Range("N134").Value = -29484.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 168750
$ws.Range("J37").Value = 168750
$ws.Range("L37").Value = 506250
$ws.Range("N37").Value = -506474

$ws.Range("H46").Value = 3687.625
$ws.Range("I46").Value = 1001
$ws.Range("J46").Value = 4583.1665
$ws.Range("K46").Value = 3003
$ws.Range("L46").Value = 13749.4995
$ws.Range("M46").Value = -2912
$ws.Range("N46").Value = -13931.4995

$ws.Range("H47").Value = 33544
$ws.Range("I47").Value = 33544
$ws.Range("K47").Value = 100632
$ws.Range("M47").Value = -100201

$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 20000
$ws.Range("K88").Value = 60000
$ws.Range("M88").Value = -59572

$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 20000
$ws.Range("K91").Value = 60000
$ws.Range("M91").Value = -58518

$ws.Range("H129").Value = 4169236
$ws.Range("I129").Value = 573.7
$ws.Range("J129").Value = 8337898.5
$ws.Range("K129").Value = 1721.1
$ws.Range("L129").Value = 25013695.5
$ws.Range("M129").Value = 3278.9
$ws.Range("N129").Value = -25023695.5

$ws.Range("H132").Value = 3262.7144
$ws.Range("J132").Value = 3657.7
$ws.Range("L132").Value = 32919.3
$ws.Range("N132").Value = -37979.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 8844.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 8844.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8844.25
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9180.25

$ws.Range("H107").Value = 1580.25
$ws.Range("I107").Value = 377.42856
$ws.Range("K107").Value = 377.42856
$ws.Range("M107").Value = 1542.57144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8831
$ws.Range("I136").Value = 3162.2307
$ws.Range("K136").Value = 9486.6921
$ws.Range("M136").Value = -6936.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 59733.75
$ws.Range("J119").Value = 59733.75
$ws.Range("L119").Value = 59733.75
$ws.Range("N119").Value = -69409.75

$ws.Range("H132").Value = 5486.9287
$ws.Range("I132").Value = 2287.5715
$ws.Range("J132").Value = 11885.643
$ws.Range("K132").Value = 6862.7145
$ws.Range("L132").Value = 35656.929
$ws.Range("M132").Value = -4332.7145
$ws.Range("N132").Value = -40716.929
